# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Swap the contents of row 16 and row 18 in the worker table:
#   - Row 16 previously held 32909336 / LEONOR MARGARITA PALOMARES LOPEZ / 2203
#   - Row 18 previously held 73193283 / TOMAS JAVIER POLO MEDRANO / 1905
# After the edit, row 16 holds the TOMAS JAVIER POLO MEDRANO record and
# row 18 holds the LEONOR MARGARITA PALOMARES LOPEZ (period 2203) record.
# Row 17 (LEONOR MARGARITA PALOMARES LOPEZ, period 2202) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> TOMAS JAVIER POLO MEDRANO record (previously on row 18)
$ws.Range("C16").Value = "73193283"
$ws.Range("D16").Value = "TOMAS JAVIER POLO MEDRANO"
$ws.Range("E16").Value = "1905"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 828116

# Row 18 -> LEONOR MARGARITA PALOMARES LOPEZ record, period 2203 (previously on row 16)
$ws.Range("C18").Value = "32909336"
$ws.Range("D18").Value = "LEONOR MARGARITA PALOMARES LOPEZ"
$ws.Range("E18").Value = "2203"
$ws.Range("F18").Value = 6933
$ws.Range("G18").Value = 1300000
